# Add the new "Shezwan House 003" sheet after the last existing sheet.
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Shezwan House 003"
$ws.Tab.Color = 5296274

# ---- Header row ----
$ws.Range("A1").Value = "SR NO"
$ws.Range("B1").Value = "ITEM DESCRIPTION"
$ws.Range("C1").Value = "QTY"
$ws.Range("D1").Value = "PRICE"
$ws.Range("E1").Value = "AMOUNT"

# ---- Item rows (A/C/D values; B filled bottom-up to reproduce shared-string order) ----
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8

$ws.Range("B9").Value = "INSTALLATION TESTING COMMISSIONING"
$ws.Range("B8").Value = "Cabling RJ 59 copper 3 + 1 Cat 6 with casing, cabling, laying "
$ws.Range("B7").Value = "Enclosure with mounting"
$ws.Range("B6").Value = "W Box 2Amps adaptor"
$ws.Range("B5").Value = "Power Supply"
$ws.Range("B4").Value = "Power Connector"
$ws.Range("B3").Value = "BNC Connector"
$ws.Range("B2").Value = "W Box  2MP Dome Camera"

$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 2000
$ws.Range("C3").Value = 12
$ws.Range("D3").Value = 60
$ws.Range("C4").Value = 6
$ws.Range("D4").Value = 40
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 1990
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 1800
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 100
$ws.Range("C8").Value = 87
$ws.Range("D8").Value = 90
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 4000

$ws.Range("E2").Formula = "=C2*D2"
$ws.Range("E3").Formula = "=C3*D3"
$ws.Range("E4").Formula = "=C4*D4"
$ws.Range("E5").Formula = "=C5*D5"
$ws.Range("E6").Formula = "=C6*D6"
$ws.Range("E7").Formula = "=C7*D7"
$ws.Range("E8").Formula = "=C8*D8"
$ws.Range("E9").Formula = "=C9*D9"

# ---- Total row ----
$ws.Range("A10").Value = "TOTAL"
$ws.Range("A10:D10").Merge()
$ws.Range("E10").Formula = "=SUM(E2:E9)"

# ---- Row heights ----
$ws.Rows.Item(8).RowHeight = 26.4

# ---- Column widths ----
$ws.Columns.Item(2).ColumnWidth = 48.44

# ---- Fonts / alignment ----
# Header row: bold, size 10, "Calibri  "
$hdr = $ws.Range("A1:E1")
$hdr.Font.Bold = $true
$hdr.Font.Size = 10
$hdr.Font.Name = "Calibri  "
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4108
$hdr.WrapText = $true

# Body rows: regular, size 10, "Calibri  "
$body = $ws.Range("A2:E9")
$body.Font.Bold = $true
$body.Font.Size = 10
$body.Font.Name = "Calibri  "
$body.Font.Bold = $false
$body.HorizontalAlignment = -4108
$body.VerticalAlignment = -4108
$body.WrapText = $true

# Total row: bold, size 10, "Calibri  "
$tot = $ws.Range("A10:E10")
$tot.Font.Bold = $true
$tot.Font.Size = 10
$tot.Font.Name = "Calibri  "
$tot.HorizontalAlignment = -4108
$tot.VerticalAlignment = -4108
$tot.WrapText = $true

# ---- Activate the new sheet and set the selection ----
$ws.Activate()
$ws.Range("H24").Select()

# ---- Restore sheet2's selection (was A8, diff shows it moves to C22) ----
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$ws2.Range("C22").Select()
$ws.Activate()
